$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 ("Participan ID" / 0 / "Participant identifier") is removed
# entirely -- participant ID becomes compulsory and is dropped from this
# configuration sheet, so the remaining rows shift up.
$ws.Rows(2).Delete()

# The header that used to read "Field" now reads "Name".
$ws.Range("A1").Value = "Name"

# Give the header columns (Name / Default) some breathing room so the
# longer field names and default values are readable.
$ws.Columns("A").ColumnWidth = 27.333333333333332
$ws.Columns("B").ColumnWidth = 48.5

# Highlight the header row.
$ws.Range("A1:C1").Interior.Color = 49407

# Add tips (as threaded comments) explaining each header column.
$ws.Range("A1").AddCommentThreaded("Name of the field to present in the startup gui")
$ws.Range("B1").AddCommentThreaded("The default value to present in the startup gui")
$ws.Range("C1").AddCommentThreaded('The corresponding "tip" to be presented when the mouse hovers over this field in the startup gui')
